$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new column (L) was inserted after K, continuing the yearly data series
# (2013..2020 in D..K) with 2021, and duplicating every other K-column
# value/format one column to the right (K3:K11 -> L3:L11).

# --- L3 : blank cell, same "header rule" formatting as K3 (bottom border) ---
$c = $ws.Range("L3")
$c.Font.Name = "Times New Roman"
$c.Font.Size = 9
$c.VerticalAlignment = -4108
$c.Borders.Item(9).LineStyle = 1
$c.Borders.Item(9).Weight = -4138
$c.Borders.Item(9).Color = 0

# --- L4 : year header, bold + right aligned + bottom border, like K4 ---
$c = $ws.Range("L4")
$c.Font.Name = "Times New Roman"
$c.Font.Size = 9
$c.Font.Bold = $true
$c.VerticalAlignment = -4108
$c.HorizontalAlignment = -4152
$c.Borders.Item(9).LineStyle = 1
$c.Borders.Item(9).Weight = -4138
$c.Borders.Item(9).Color = 0
$c.Value = 2021

# --- L5:L10 : plain data values, same formatting as K5:K10 ---
foreach ($r in 5..10) {
    $c = $ws.Range("L$r")
    $c.Font.Name = "Times New Roman"
    $c.Font.Size = 9
    $c.VerticalAlignment = -4108
}
$ws.Range("L5").Value = 0.86
$ws.Range("L6").Value = 1.07
$ws.Range("L7").Value = 25.27
$ws.Range("L8").Value = 14
$ws.Range("L9").Value = 0.12
$ws.Range("L10").Value = 21.74

# --- L11 : totals row, bottom border like K11 ---
$c = $ws.Range("L11")
$c.Font.Name = "Times New Roman"
$c.Font.Size = 9
$c.VerticalAlignment = -4108
$c.Borders.Item(9).LineStyle = 1
$c.Borders.Item(9).Weight = -4138
$c.Borders.Item(9).Color = 0
$c.Value = 9.4600000000000009

# The active selection moved to N2.
$ws.Range("N2").Select()
